{"js": "// Update the date line and every \"A\u00f7B=C, D\" answer cell in the table to the\n// new values from the commit. Each old value is unique within the document,\n// so a direct search-and-replace (by exact text) for each pair is safe and\n// unambiguous.\nconst replacements = [\n  [\"2026-01-05 Monday\", \"2026-01-06 Tuesday\"],\n  [\"289\u00f72=144, 1\", \"802\u00f79=89, 1\"],\n  [\"110\u00f76=18, 2\", \"952\u00f77=136, 0\"],\n  [\"838\u00f77=119, 5\", \"269\u00f74=67, 1\"],\n  [\"721\u00f79=80, 1\", \"822\u00f73=274, 0\"],\n  [\"464\u00f78=58, 0\", \"348\u00f78=43, 4\"],\n  [\"742\u00f76=123, 4\", \"485\u00f76=80, 5\"],\n  [\"967\u00f73=322, 1\", \"931\u00f75=186, 1\"],\n  [\"793\u00f79=88, 1\", \"103\u00f77=14, 5\"],\n  [\"627\u00f73=209, 0\", \"583\u00f77=83, 2\"],\n  [\"145\u00f76=24, 1\", \"452\u00f75=90, 2\"],\n  [\"981\u00f78=122, 5\", \"176\u00f72=88, 0\"],\n  [\"200\u00f79=22, 2\", \"107\u00f73=35, 2\"],\n  [\"992\u00f79=110, 2\", \"997\u00f72=498, 1\"],\n  [\"612\u00f77=87, 3\", \"437\u00f72=218, 1\"],\n  [\"684\u00f75=136, 4\", \"678\u00f78=84, 6\"],\n  [\"244\u00f74=61, 0\", \"357\u00f74=89, 1\"],\n  [\"468\u00f75=93, 3\", \"464\u00f72=232, 0\"],\n  [\"914\u00f77=130, 4\", \"965\u00f75=193, 0\"],\n  [\"812\u00f73=270, 2\", \"604\u00f78=75, 4\"],\n  [\"101\u00f74=25, 1\", \"406\u00f77=58, 0\"],\n  [\"960\u00f73=320, 0\", \"580\u00f78=72, 4\"],\n  [\"192\u00f74=48, 0\", \"667\u00f72=333, 1\"],\n  [\"419\u00f77=59, 6\", \"874\u00f79=97, 1\"],\n  [\"788\u00f78=98, 4\", \"565\u00f72=282, 1\"],\n  [\"468\u00f79=52, 0\", \"808\u00f72=404, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every \"A\u00f7B=C, D\" answer cell in the table to the\n# new values from the commit. Each old value is unique within the document,\n# so a Find/Replace (one occurrence each) for every pair is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-01-05 Monday\", \"2026-01-06 Tuesday\"),\n    @(\"289\u00f72=144, 1\", \"802\u00f79=89, 1\"),\n    @(\"110\u00f76=18, 2\", \"952\u00f77=136, 0\"),\n    @(\"838\u00f77=119, 5\", \"269\u00f74=67, 1\"),\n    @(\"721\u00f79=80, 1\", \"822\u00f73=274, 0\"),\n    @(\"464\u00f78=58, 0\", \"348\u00f78=43, 4\"),\n    @(\"742\u00f76=123, 4\", \"485\u00f76=80, 5\"),\n    @(\"967\u00f73=322, 1\", \"931\u00f75=186, 1\"),\n    @(\"793\u00f79=88, 1\", \"103\u00f77=14, 5\"),\n    @(\"627\u00f73=209, 0\", \"583\u00f77=83, 2\"),\n    @(\"145\u00f76=24, 1\", \"452\u00f75=90, 2\"),\n    @(\"981\u00f78=122, 5\", \"176\u00f72=88, 0\"),\n    @(\"200\u00f79=22, 2\", \"107\u00f73=35, 2\"),\n    @(\"992\u00f79=110, 2\", \"997\u00f72=498, 1\"),\n    @(\"612\u00f77=87, 3\", \"437\u00f72=218, 1\"),\n    @(\"684\u00f75=136, 4\", \"678\u00f78=84, 6\"),\n    @(\"244\u00f74=61, 0\", \"357\u00f74=89, 1\"),\n    @(\"468\u00f75=93, 3\", \"464\u00f72=232, 0\"),\n    @(\"914\u00f77=130, 4\", \"965\u00f75=193, 0\"),\n    @(\"812\u00f73=270, 2\", \"604\u00f78=75, 4\"),\n    @(\"101\u00f74=25, 1\", \"406\u00f77=58, 0\"),\n    @(\"960\u00f73=320, 0\", \"580\u00f78=72, 4\"),\n    @(\"192\u00f74=48, 0\", \"667\u00f72=333, 1\"),\n    @(\"419\u00f77=59, 6\", \"874\u00f79=97, 1\"),\n    @(\"788\u00f78=98, 4\", \"565\u00f72=282, 1\"),\n    @(\"468\u00f79=52, 0\", \"808\u00f72=404, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $null = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $true, $newText, 2)\n}\n"}
